# Remove the "Squad Total" summary row (row 35): clear the label in A35
# and all the SUM/AVERAGE formulas across B35:U35 so the row becomes
# entirely empty (only the existing cell styles remain).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A35:U35").ClearContents() | Out-Null

# Update the active selection to match the saved view (R2).
$ws.Range("R2").Select() | Out-Null
